$d = $word.ActiveDocument

# The whole body is being collapsed down to a single paragraph. Start by
# merging every paragraph after the first one into paragraph 1 (deleting
# from the end of paragraph 1's text through the end of the document
# removes all the intervening paragraph marks along with their content).
$p1 = $d.Paragraphs.Item(1)
$tail = $d.Range($p1.Range.End, $d.Content.End)
$tail.Delete()

# Insert the new paragraph (with its proofErr markers and the relocated
# _GoBack bookmark) immediately in front of the surviving paragraph.
$insertPoint = $d.Range(0, 0)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">The data is </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>actually normally</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> distributed, but it might need transformation to reveal its normality. For example, lognormal distribution </w:t></w:r>
<w:r><w:t>bec</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t xml:space="preserve">omes normal distribution after taking a log on it. </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertPoint.InsertXML($xml)

# Remove the now-redundant leftover paragraph (the document's original
# first paragraph, "Motivation") that was pushed after our inserted text.
$newFirst = $d.Paragraphs.Item(1)
$leftover = $d.Range($newFirst.Range.End, $d.Content.End)
$leftover.Delete()
